# Add flowthrough data to separate dataframe and update output file
#
# Appends two new sample groups ("flowthrough1" and "flowthrough1_no_beads")
# below the existing enrich1/enrich2 data (rows 2-19), growing the used
# range from A1:E19 to A1:E37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @(20, 0,  "flowthrough1",          "RNA18S1",  "-23.511945",         "8.359865661106497e-06"),
  @(21, 1,  "flowthrough1",          "RNA18S1",  "-21.839316",         "2.665083369387191e-05"),
  @(22, 2,  "flowthrough1",          "RNA18S1",  "-22.082388",         "2.251846775642507e-05"),
  @(23, 3,  "flowthrough1",          "vtRNA1-1", "-9.168227000000002", "0.1738155655915517"),
  @(24, 4,  "flowthrough1",          "vtRNA1-1", "-8.392045999999997", "0.2976750457239626"),
  @(25, 5,  "flowthrough1",          "vtRNA1-1", "-8.591697",          "0.2592038747363474"),
  @(26, 6,  "flowthrough1",          "FFLUC",    "-4.581559000000002", "4.176507966076688"),
  @(27, 7,  "flowthrough1",          "FFLUC",    "-5.376252000000001", "2.407614131490228"),
  @(28, 8,  "flowthrough1",          "FFLUC",    "-5.448609999999999", "2.289839039155415"),
  @(29, 9,  "flowthrough1_no_beads", "RNA18S1",  $null, $null),
  @(30, 10, "flowthrough1_no_beads", "RNA18S1",  $null, $null),
  @(31, 11, "flowthrough1_no_beads", "RNA18S1",  $null, $null),
  @(32, 12, "flowthrough1_no_beads", "vtRNA1-1", $null, $null),
  @(33, 13, "flowthrough1_no_beads", "vtRNA1-1", $null, $null),
  @(34, 14, "flowthrough1_no_beads", "vtRNA1-1", $null, $null),
  @(35, 15, "flowthrough1_no_beads", "FFLUC",    $null, $null),
  @(36, 16, "flowthrough1_no_beads", "FFLUC",    $null, $null),
  @(37, 17, "flowthrough1_no_beads", "FFLUC",    $null, $null)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    if ($row[4] -ne $null) {
        $ws.Cells.Item($r, 4).Value = [double]$row[4]
    }
    if ($row[5] -ne $null) {
        $ws.Cells.Item($r, 5).Value = [double]$row[5]
    }
}

# Column A in the existing data carries the bold/centered/bordered style
# (index 1 in styles.xml) - copy that formatting down onto the new rows.
$ws.Range("A2").Copy()
$ws.Range("A20:A37").PasteSpecial(-4122)  # xlPasteFormats

# The "no_beads" rows (29-37) have no CT.difference / percent_recovery
# measurements, but the source data still emits empty cells in D/E for
# them (matching a pandas export of NaN). Copy the plain (unstyled)
# format from an existing data cell so the cells exist but stay blank.
$ws.Range("D19").Copy()
$ws.Range("D29:D37").PasteSpecial(-4122)
$ws.Range("E19").Copy()
$ws.Range("E29:E37").PasteSpecial(-4122)

Write-Output "Added flowthrough1 / flowthrough1_no_beads rows (20-37)."
